$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TruthFulQA")

# Insert a new row for the "L11" model result, right after "L10_last"
# (row 14) and before "Q1" (old row 15). This pushes the existing
# Q1..Q7 rows down by one.
$ws.Rows("15:15").Insert()

$ws.Range("A15").Value = "L11"
$ws.Range("B15").Value = 0.41860465116279
$ws.Range("C15").Value = 0.43084455324357401
$ws.Range("D15").Value = 0.476132190942472
$ws.Range("E15").Formula = "=AVERAGE(Table3[[#This Row],[bleu_acc]:[bluert_acc]])"

# The TruthFulQA sheet's data lives inside Table3 (A1:E21) - grow it by
# one row now that the new "L11" row has been added.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E22"))

# Match the author's final UI state: TruthFulQA tab active, selection
# resting near the newly added row.
$ws.Activate() | Out-Null
$ws.Range("I15").Select() | Out-Null
